# Updated NGIN test data sheet: bump the "1020" test-data instance to "1021"
# across every data cell in row 2 that encodes the instance number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A2"  = "NGIN1021"
    "C2"  = "NGIN1021"
    "D2"  = "ngindomain1021.com"
    "F2"  = "nginocn1021"
    "G2"  = "testreference1021"
    "H2"  = "ngincontact1021"
    "J2"  = "ngin1021@test.com"
    "Q2"  = "NGIN1021"
    "T2"  = "NGINOrder_1021"
    "U2"  = "NGINRFI_1021"
    "V2"  = "NGINOrder_1021"
    "W2"  = "NGINRFI_1021"
    "Y2"  = "NGINService_1021"
    "AJ2" = "NGINUser_1021"
    "AK2" = "User_1021"
    "AN2" = "NGINUser_1021@gmail.com"
    "AP2" = "NGINOrder_1021"
    "AQ2" = "NGINUseredit1021"
    "AR2" = "Useredit1021"
    "AU2" = "NGINUseredit_1021@gmail.com"
    "AZ2" = "NGINOrderedit_1021"
    "BA2" = "NGINRFIedit_1021"
    "BB2" = "NGINOrder_1021"
    "BC2" = "NGINRFI_1021"
    "BI2" = "nginreseller1021@gmail.com"
    "BQ2" = "nginreselleredit1021@gmail.com"
    "BY2" = "Reseller1021"
    "DP2" = "Reselleredit1021"
    "ED2" = "AT-nginocn1021"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# EF2 holds a numeric-looking value stored as text (quote-prefixed in the
# original file) - keep the leading apostrophe so it stays text instead of
# being reinterpreted as a number.
$ws.Range("EF2").Value = "'390201021891"

# Sheet view now scrolled back to the top-left with C3 selected (was topLeftCell DZ1 / EE6 selected)
$ws.Range("C3").Select()
